$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("imc")

$ws.Range("B86").Value = 1.5
$ws.Range("B87").Value = 1.4
